$wb = $excel.ActiveWorkbook

# --- Add a new "Player Info" sheet, placed before the existing "ODI Batting" sheet ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Player Info"

# Header row
$newSheet.Range("A1").Value = "ID"
$newSheet.Range("B1").Value = "NAME"
$newSheet.Range("C1").Value = "BATTING_HAND"
$newSheet.Range("D1").Value = "BOWL_STYLE"

# Data row - force column A to text so the numeric-looking ID stays a string
$newSheet.Range("A2").NumberFormat = "@"
$newSheet.Range("A2").Value = "5956"
$newSheet.Range("B2").Value = "Finnley Hugh Allen"
$newSheet.Range("C2").Value = "Right Handed"
$newSheet.Range("D2").Value = "Does Not Bowl | Unknown"
$newSheet.Range("A2").Style = "Normal"

# Match the header styling used on the other sheet (bold, bordered, centered/top-aligned)
$headerRange = $newSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# --- Update the "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$ws = $wb.Worksheets.Item("ODI Batting")
$ws.Range("D1").Value = "MATCH_CODE"

$lastRow = $ws.UsedRange.Rows.Count
$dataRange = $ws.Range("D2:D" + $lastRow)
$dataRange.NumberFormat = "@"
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $text = $cell.Text
    if ($text -match "MatchCode=(\d+)") {
        $code = $matches[1]
        $cell.Value = $code
    }
}
$dataRange.Style = "Normal"

Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
